# Estadisticos Segundo Parcial 23 Mayo
#
# 1. "Estadisticos 2P"   -> fill in the second-partial results (D:H) that were
#                            previously placeholders (all-zero columns).
# 2. "Estadisticos Final" -> refresh E:H with the newly combined totals.
# 3. "Rescatables"        -> regenerate the make-up-work list: it grows from
#                            11 students to 17, re-sorted by # of failed
#                            subjects (column G) descending.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Estadisticos 2P"
# ---------------------------------------------------------------------------
$ws2p = $wb.Worksheets.Item("Estadisticos 2P")

$ws2p.Cells.Item(2, 4).Value = 0
$ws2p.Cells.Item(2, 5).Value = 9
$ws2p.Cells.Item(2, 6).Value = 22
$ws2p.Cells.Item(2, 7).Value = 70.97
$ws2p.Cells.Item(2, 8).Value = 6.4

$ws2p.Cells.Item(3, 4).Value = 0
$ws2p.Cells.Item(3, 5).Value = 11
$ws2p.Cells.Item(3, 6).Value = 38
$ws2p.Cells.Item(3, 7).Value = 77.55
$ws2p.Cells.Item(3, 8).Value = 6.3

$ws2p.Cells.Item(4, 4).Value = 0
$ws2p.Cells.Item(4, 5).Value = 0
$ws2p.Cells.Item(4, 6).Value = 31
$ws2p.Cells.Item(4, 7).Value = 100
$ws2p.Cells.Item(4, 8).Value = 7.3

$ws2p.Cells.Item(5, 4).Value = 0
$ws2p.Cells.Item(5, 5).Value = 2
$ws2p.Cells.Item(5, 6).Value = 33
$ws2p.Cells.Item(5, 7).Value = 94.29
$ws2p.Cells.Item(5, 8).Value = 7.9

# ---------------------------------------------------------------------------
# Sheet "Estadisticos Final"
# ---------------------------------------------------------------------------
$wsFinal = $wb.Worksheets.Item("Estadisticos Final")

$wsFinal.Cells.Item(2, 5).Value = 9
$wsFinal.Cells.Item(2, 6).Value = 22
$wsFinal.Cells.Item(2, 7).Value = 70.97
$wsFinal.Cells.Item(2, 8).Value = 6.7

$wsFinal.Cells.Item(3, 8).Value = 7

$wsFinal.Cells.Item(4, 5).Value = 0
$wsFinal.Cells.Item(4, 6).Value = 31
$wsFinal.Cells.Item(4, 7).Value = 100
$wsFinal.Cells.Item(4, 8).Value = 7.7

$wsFinal.Cells.Item(5, 5).Value = 2
$wsFinal.Cells.Item(5, 6).Value = 33
$wsFinal.Cells.Item(5, 7).Value = 94.29
$wsFinal.Cells.Item(5, 8).Value = 7.8

# ---------------------------------------------------------------------------
# Sheet "Rescatables"
# ---------------------------------------------------------------------------
$wsResc = $wb.Worksheets.Item("Rescatables")

# Wipe the old 11-student table; it gets fully replaced below.
$wsResc.Range("A2:G12").Clear()

$materiaCD = "Cultura digital II"
$materiaDAM = "DESARROLLA APLICACIONES M`u{00D3}VILES PARA IOS"

$data = @(
    @(24330051920259, "RUIZ",      "PELLICO",   "YOSHUA RAFAEL",       $materiaCD,  "2APM",  4),
    @(24330051920117, "PLIEGO",    "TORRES",    "MYA YAMILET",         $materiaCD,  "2ARHM", 4),
    @(24330051920348, "TRUJILLO",  "FLORES",    "VALERIA",             $materiaCD,  "2ARHM", 4),
    @(24330051920387, "NICANOR",   "MALDONADO", "DENISSE ARELI",       $materiaCD,  "2APM",  3),
    @(24330051920129, "DIAZ",      "MARCELINO", "JUAN",                $materiaCD,  "2ARHM", 3),
    @(24330051920345, "GARIBAY",   "GOMEZ",     "LIZBETH MARIAM",      $materiaCD,  "2ARHM", 3),
    @(22330051920426, "VICENTE",   "QUINTANA",  "CRISTIAN",            $materiaDAM, "6APM",  3),
    @(24330051920384, "RAMIREZ",   "BELLO",     "JOACIM ALBERTO",      $materiaCD,  "2APM",  2),
    @(24330051920350, "FLORES",    "LOBATO",    "MARIANA",             $materiaCD,  "2ARHM", 2),
    @(24330051920132, "GONZALEZ",  "CRUZ",      "JESUS",               $materiaCD,  "2ARHM", 2),
    @(24330051920351, "PLIEGO",    "LORENZO",   "CALEB SANTIAGO",      $materiaCD,  "2ARHM", 2),
    @(24330051920279, "GUTIERREZ", "HUERTA",    "DIEGO",               $materiaCD,  "2APM",  1),
    @(24330051920372, "LOPEZ",     "GONZALEZ",  "GUADALUPE",           $materiaCD,  "2APM",  1),
    @(24330051920250, "ESTEVEZ",   "MARIN",     "ESTRELLA MONTSERRAT", $materiaCD,  "2ARHM", 1),
    @(24330051920281, "MONTERD",   "GARCIA",    "ISAI",                $materiaCD,  "2ARHM", 1),
    @(24330051920347, "MOLINA",    "MACUISTLE", "ANGEL FRANCISCO",     $materiaCD,  "2ARHM", 1),
    @(22330051920172, "ROJAS",     "ANGUIANO",  "LUIS ANGEL",          $materiaDAM, "6APM",  1)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $rec = $data[$i]
    $wsResc.Cells.Item($row, 1).Value = $rec[0]
    $wsResc.Cells.Item($row, 2).Value = $rec[1]
    $wsResc.Cells.Item($row, 3).Value = $rec[2]
    $wsResc.Cells.Item($row, 4).Value = $rec[3]
    $wsResc.Cells.Item($row, 5).Value = $rec[4]
    $wsResc.Cells.Item($row, 6).Value = $rec[5]
    $wsResc.Cells.Item($row, 7).Value = $rec[6]
}
